# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The value 45186 (serial date) is updated to 45188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
